$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are numeric-looking strings ("29.328.24",
# "0.9987", "7.400", ...) that must stay TEXT, matching the source
# t="inlineStr" cells. A bare assignment lets Excel auto-convert them to
# real numbers (losing formatting like trailing zeros / multi-dot groups),
# so each is written with a leading apostrophe (forces text entry) and
# then ClearFormats() removes the transient "quote prefix" cell style so
# the cell keeps using the same (default) style as before the edit.

$ws.Range("D2").Value = "'29.328.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "'1.842.12"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "'0.9987"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").Value = "'238.95"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'0.6291"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.07515"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "'0.2934"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'24.40"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'0.07693"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "'1.833.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.65%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'0.6775"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "'0.00001046"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.25%  "
$ws.Range("D16").Value = "'82.87"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'2.071.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -8.54%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Value = "'29.359.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "'227.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'0.9997"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'7.400"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.64%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'156.49"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "'0.1386"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "'8.347"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "'17.58"
$ws.Range("D28").ClearFormats()
$ws.Range("D29").Value = "'1.454"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "'0.05618"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "'4.096"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'4.015"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "'1.154"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "'0.7078"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").Value = "'2.590"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'1.238.59"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("D39").Value = "'0.01808"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'2.760"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").Value = "'6.241"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").Value = "'0.9003"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").Value = "'0.9992"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'101.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "'65.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("D46").Value = "'0.00000000120"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "'0.3986"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("D49").Value = "'8.915"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").Value = "'1.666"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -0.34%  "
